$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Flip the sign of the EmployeeId values for the first three rows
$ws.Range("D2").Value = -1
$ws.Range("D3").Value = -1
$ws.Range("D4").Value = -1

# Update the active selection/view: active cell D3, no scrolled topLeftCell
$ws.Range("D3").Select() | Out-Null
